$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title (B3) and description (C3) for the Dr. Kunhao Yang news item
$ws.Range("B3").Value = "Dr. Kunhao Yang Joined SIT College of Engineering"
$ws.Range("C3").Value = " Dr. Kunhao Yang has joined Shibaura Institute of Technology as an Assistant Professor in the College of Engineering and as the director of the Computational Social Science (CSS) Laboratory. Dr. Yang will be spearheading research at the intersection of social science and data science, and will also be teaching new courses exploring these dynamic fields."

# Recalculate row height now that the text is shorter
$ws.Rows.Item(3).RowHeight = 102

# Update the selected cell/active view to C3 (was F3)
$ws.Range("C3").Select()
